# Update Excel files from OneDrive - Wed May  7 18:46:18 UTC 2025
$wb = $excel.ActiveWorkbook

# The "COMPLETED" sheet is the active/selected sheet in this workbook.
$ws = $wb.Worksheets.Item("COMPLETED")
$ws.Activate()

# Row 2 was a placeholder/"in progress" entry (PRODUCT/CERTIFICATION/AGENCY all
# showed the filler text "f"). Fill in the real values for this completed
# certification record.
$ws.Range("A2").Value = "PARADEA"
$ws.Range("B2").Value = "IEC 61215"
$ws.Range("C2").Value = "TUV"

# The NOTES column filler text "f" is now replaced with "OK" wherever it still
# appears (row 2 and row 3).
$ws.Range("E2").Value = "OK"
$ws.Range("E3").Value = "OK"

# Move the active cell selection on the sheet down one row (F7 -> F8).
$ws.Range("F8").Select()

# Scroll the workbook tabs strip over by one sheet.
$win = $excel.ActiveWindow
$win.ScrollWorkbookTabs(1)
